$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.313.62"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.843.15"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'240.07"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'0.6278"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'0.07444"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").Value = "'0.2897"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "'0.07735"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.843.67"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").Value = "'4.977"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "'0.6789"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "'0.00001043"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "'81.86"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "'6.173"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "29.367.57"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'227.75"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").Value = "'12.29"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'7.494"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'1.0000"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "'159.29"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'8.479"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").Value = "'0.1368"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").Value = "'17.48"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "'0.06514"
$ws.Range("E28").Value = "  +15.99%  "
$ws.Range("D29").Value = "'1.424"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").Value = "'1.481"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Value = "'4.086"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "'4.081"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "'1.831"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "'1.139"
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("D35").Value = "'0.6938"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "'2.578"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "1.259.31"
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("D38").Value = "'2.831"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").Value = "'0.01831"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").Value = "'6.717"
$ws.Range("E40").Value = "  +4.56%  "
$ws.Range("D41").Value = "'0.9252"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "'0.9989"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "2.006.36"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "'101.26"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'65.84"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000119"
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.724"
$ws.Range("E47").Value = "  +2.30%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.048"
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1152"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.986"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.3924"
$ws.Range("E51").Value = "  -2.08%  "
